$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new blank rows at row 847, pushing the existing rows 847-860
# down to 854-867 (unchanged content).
$ws.Rows("847:853").Insert()

# New weekly price rows (Palta / Comercializadora del Agro de Limari),
# all sharing the same constant columns as the surrounding data.
$newRows = @(
    @{ Row=847; D=45239; K='Edranol'; L='Especial'; M=240; N=2200; O=2300; P=2250; Q='$/kilo (en caja de 17 kilos)'; S=2250 },
    @{ Row=848; D=45239; K='Edranol'; L='Primera';  M=300; N=1900; O=2000; P=1950; Q='$/kilo (en caja de 17 kilos)'; S=1950 },
    @{ Row=849; D=45239; K='Edranol'; L='Segunda';  M=240; N=1600; O=1700; P=1650; Q='$/kilo (en caja de 17 kilos)'; S=1650 },
    @{ Row=850; D=45239; K='Hass';    L='Especial'; M=240; N=2800; O=2900; P=2850; Q='$/kilo (en caja de 17 kilos)'; S=2850 },
    @{ Row=851; D=45239; K='Hass';    L='Primera';  M=400; N=2600; O=2700; P=2650; Q='$/kilo (en caja de 17 kilos)'; S=2650 },
    @{ Row=852; D=45239; K='Hass';    L='Segunda';  M=360; N=2100; O=2200; P=2150; Q='$/kilo (en caja de 17 kilos)'; S=2150 },
    @{ Row=853; D=45239; K='Hass';    L='Tercera';  M=200; N=1800; O=1900; P=1850; Q='$/kilo (en caja de 17 kilos)'; S=1850 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2 = 2
    $ws.Cells.Item($row, 2).Value2 = 'Comercializadora del Agro de Limarí'
    $ws.Cells.Item($row, 3).Value2 = 'Coquimbo'
    $ws.Cells.Item($row, 4).Value2 = $r.D
    $ws.Cells.Item($row, 5).Value2 = 4
    $ws.Cells.Item($row, 6).Value2 = 'Fruta'
    $ws.Cells.Item($row, 7).Value2 = 100106
    $ws.Cells.Item($row, 8).Value2 = 'Oleaginosos'
    $ws.Cells.Item($row, 9).Value2 = 100106002
    $ws.Cells.Item($row, 10).Value2 = 'Palta'
    $ws.Cells.Item($row, 11).Value2 = $r.K
    $ws.Cells.Item($row, 12).Value2 = $r.L
    $ws.Cells.Item($row, 13).Value2 = $r.M
    $ws.Cells.Item($row, 14).Value2 = $r.N
    $ws.Cells.Item($row, 15).Value2 = $r.O
    $ws.Cells.Item($row, 16).Value2 = $r.P
    $ws.Cells.Item($row, 17).Value2 = $r.Q
    $ws.Cells.Item($row, 18).Value2 = 'Provincia de Limarí'
    $ws.Cells.Item($row, 19).Value2 = $r.S
    $ws.Cells.Item($row, 20).Value2 = 1
}
